# Add 2022-Q4 data:
#  - Insert a brand-new "2022-Q4" worksheet right after "总计", holding the
#    same shape/headers as the existing quarterly sheets, populated with the
#    new Q4 fund holdings.
#  - Update the "总计" (summary) worksheet: a new row for 2022-Q4 is inserted
#    at the top of the data (row 2) and the previously existing rows shift
#    down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet ("总计"): shift existing rows down and insert the new
#    2022-Q4 figures at the top of the table.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Extend the index-column style (s="2") down to the newly-added row 7.
$total.Range("A2").Copy()
$total.Range("A7").PasteSpecial(-4122)

$summaryRows = @(
    @(0, "2022-Q4", 6,  0.34),
    @(1, "2022-Q3", 6,  0.04),
    @(2, "2022-Q2", 10, 0.76),
    @(3, "2022-Q1", 7,  1.7),
    @(4, "2021-Q4", 6,  0.93),
    @(5, "2021-Q3", 22, 5.53)
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = 2 + $i
    $total.Cells.Item($r, 1).Value = $summaryRows[$i][0]
    $total.Cells.Item($r, 2).Value = $summaryRows[$i][1]
    $total.Cells.Item($r, 3).Value = $summaryRows[$i][2]
    $total.Cells.Item($r, 4).Value = $summaryRows[$i][3]
}

# ---------------------------------------------------------------------
# 2) New "2022-Q4" sheet: copy the layout/styling of the existing
#    "2022-Q3" sheet (same header row + formatting), place it right after
#    "总计", rename it, then overwrite the data with the Q4 figures.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $total)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Columns B (fund code) and D:G (numeric-looking text figures) must stay
# text so leading zeros / exact string formatting survive, matching the
# rest of the workbook's convention for these columns.
$q4.Range("B2:B7").NumberFormat = "@"
$q4.Range("D2:G7").NumberFormat = "@"

$fundRows = @(
    @("005870", "鹏华沪深300指数增强A", "11.31", "92.77", "2.07",  "0.2341",  9),
    @("016690", "鹏华沪深300指数增强C", "3.70",  "92.77", "2.07",  "0.0766",  9),
    @("014938", "同泰产业升级混合A",    "0.95",  "68.98", "2.74",  "0.0260",  10),
    @("080007", "长盛同鑫行业配置混合A", "0.20",  "84.46", "2.30",  "0.0046",  9),
    @("010991", "长盛同鑫行业配置混合C", "0.02",  "84.46", "2.30",  "0.0005",  9),
    @("014939", "同泰产业升级混合C",    "-0.01", "68.98", "2.74",  "-0.0003", 10)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = 2 + $i
    $q4.Cells.Item($r, 2).Value = $fundRows[$i][0]
    $q4.Cells.Item($r, 3).Value = $fundRows[$i][1]
    $q4.Cells.Item($r, 4).Value = $fundRows[$i][2]
    $q4.Cells.Item($r, 5).Value = $fundRows[$i][3]
    $q4.Cells.Item($r, 6).Value = $fundRows[$i][4]
    $q4.Cells.Item($r, 7).Value = $fundRows[$i][5]
    $q4.Cells.Item($r, 8).Value = $fundRows[$i][6]
}
